$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 (Sending=FAPs, Target=FAPs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 4.278141666666667
$ws.Range("H2").Value = 12.834425
$ws.Range("I2").Value = 0.9663225094340192
$ws.Range("J2").Value = 0.9663225094340191
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.292736333333333
$ws.Range("N2").Value = 3.878209
$ws.Range("O2").Value = 0.9024289107934966
$ws.Range("P2").Value = 0.9024289107934966
$ws.Range("Q2").Value = 5.530509171647222
$ws.Range("R2").Value = 49.774582544825
$ws.Range("S2").Value = 0.8720373696637803
$ws.Range("T2").Value = 0.8720373696637802

# Update existing row 3, now becomes Sending=FAPs, Target=sCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 4.278141666666667
$ws.Range("H3").Value = 12.834425
$ws.Range("I3").Value = 0.9663225094340192
$ws.Range("J3").Value = 0.9663225094340191
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.1397713333333333
$ws.Range("N3").Value = 0.419314
$ws.Range("O3").Value = 0.09757108920650338
$ws.Range("P3").Value = 0.09757108920650338
$ws.Range("Q3").Value = 0.5979615649388889
$ws.Range("R3").Value = 5.38165408445
$ws.Range("S3").Value = 0.09428513977023889
$ws.Range("T3").Value = 0.09428513977023889

# New row 4: Sending=sCs, Target=FAPs
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Bmp7"
$ws.Range("C4").Value = "Bmpr1b"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.1490983333333333
$ws.Range("H4").Value = 0.447295
$ws.Range("I4").Value = 0.03367749056598092
$ws.Range("J4").Value = 0.03367749056598091
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.292736333333333
$ws.Range("N4").Value = 3.878209
$ws.Range("O4").Value = 0.9024289107934966
$ws.Range("P4").Value = 0.9024289107934966
$ws.Range("Q4").Value = 0.1927448327394444
$ws.Range("R4").Value = 1.734703494655
$ws.Range("S4").Value = 0.03039154112971642
$ws.Range("T4").Value = 0.03039154112971641

# New row 5: Sending=sCs, Target=sCs
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Bmp7"
$ws.Range("C5").Value = "Bmpr1b"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.1490983333333333
$ws.Range("H5").Value = 0.447295
$ws.Range("I5").Value = 0.03367749056598092
$ws.Range("J5").Value = 0.03367749056598091
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.1397713333333333
$ws.Range("N5").Value = 0.419314
$ws.Range("O5").Value = 0.09757108920650338
$ws.Range("P5").Value = 0.09757108920650338
$ws.Range("Q5").Value = 0.02083967284777778
$ws.Range("R5").Value = 0.18755705563
$ws.Range("S5").Value = 0.0032859494362645
$ws.Range("T5").Value = 0.0032859494362645
